$d = $word.ActiveDocument

# 1) "Total de citas programadas: " run stays untouched; only the "4" run becomes "7".
$para = $d.Paragraphs.Item(3)
$para.Range.Find.Execute("4", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "7", 2) | Out-Null

# 2) Work on the first (only) table.
$t = $d.Tables.Item(1)

# Row 2 (08:30 - 08:45): INTERLINK2AMERICAS -> PROCOLOMBIA
$t.Cell(2, 3).Range.Text = "PROCOLOMBIA"

# Insert new row after row 2: 08:45 - 09:00 | (empty) | ARMANDO VELÁSQUEZ
$newRow = $t.Rows.Add($t.Rows.Item(3))
$t.Cell(3, 1).Range.Text = "08:45 - 09:00"
$t.Cell(3, 3).Range.Text = "ARMANDO VELÁSQUEZ"

# Insert new row after that: 09:00 - 09:15 | (empty) | REGIONAL S.A.S
$newRow = $t.Rows.Add($t.Rows.Item(4))
$t.Cell(4, 1).Range.Text = "09:00 - 09:15"
$t.Cell(4, 3).Range.Text = "REGIONAL S.A.S"

# Row 5 (09:15 - 09:30 / BOX BRAND) remains unchanged.

# Insert new row before the "09:45 - 10:00" row (currently row 6):
# 09:30 - 09:45 | (empty) | CAFÉ MOLINA
$newRow = $t.Rows.Add($t.Rows.Item(6))
$t.Cell(6, 1).Range.Text = "09:30 - 09:45"
$t.Cell(6, 3).Range.Text = "CAFÉ MOLINA"

# Row 7 (09:45 - 10:00): ARMANDO VELÁSQUEZ -> FLOR A FRUTO
$t.Cell(7, 3).Range.Text = "FLOR A FRUTO"

# Row 8 (10:00 - 10:15): FLOR A FRUTO -> INTERLINK2AMERICAS
$t.Cell(8, 3).Range.Text = "INTERLINK2AMERICAS"

Write-Output ("Final row count: " + $t.Rows.Count)
